# Update the cryptos list (Price and Volume(1h) columns) with the latest
# scraped values, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that look numeric (e.g. "0.9994") must be forced to Text
# format first, otherwise Excel auto-converts them to numbers and the
# original formatting (e.g. trailing zeros, dotted thousand separators) is lost.
$textPriceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D16",
    "D18",
    "D19",
    "D20",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D35",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D48",
    "D50",
    "D51",
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    2 = @{ D = "29.261.50"; E = "  -0.41%  " }
    3 = @{ D = "1.861.84"; E = "  -1.10%  " }
    4 = @{ D = "0.9994"; E = "  -0.11%  " }
    5 = @{ D = "0.7051"; E = "  -1.23%  " }
    6 = @{ D = "242.51"; E = "  -0.20%  " }
    7 = @{ D = "0.9995"; E = "  -0.10%  " }
    8 = @{ D = "0.3146"; E = "  +0.48%  " }
    9 = @{ D = "0.07811"; E = "  -2.89%  " }
    10 = @{ D = "24.29"; E = "  -3.69%  " }
    11 = @{ D = "0.08005"; E = "  -4.11%  " }
    12 = @{ D = "1.861.96"; E = "  -1.23%  " }
    13 = @{ D = "94.10"; E = "  -0.35%  " }
    14 = @{ D = "5.186"; E = "  -1.36%  " }
    15 = @{ E = "  -3.04%  " }
    16 = @{ D = "6.417"; E = "  +1.39%  " }
    17 = @{ D = "29.259.00"; E = "  -0.44%  " }
    18 = @{ D = "0.000008285"; E = "  -3.07%  " }
    19 = @{ D = "253.57"; E = "  +4.67%  " }
    20 = @{ D = "13.14"; E = "  -0.96%  " }
    21 = @{ D = "2.107.90"; E = "  -1.39%  " }
    22 = @{ E = "  -0.14%  " }
    23 = @{ D = "7.564"; E = "  -3.92%  " }
    24 = @{ D = "0.9994"; E = "  -0.14%  " }
    25 = @{ D = "0.1565"; E = "  -1.47%  " }
    26 = @{ E = "  -0.96%  " }
    27 = @{ D = "160.08"; E = "  -2.06%  " }
    28 = @{ D = "18.89"; E = "  +1.29%  " }
    29 = @{ D = "1.493"; E = "  -1.28%  " }
    30 = @{ D = "4.312"; E = "  -2.54%  " }
    31 = @{ E = "  -1.41%  " }
    32 = @{ D = "1.208"; E = "  +0.63%  " }
    33 = @{ D = "0.05283"; E = "  -1.97%  " }
    34 = @{ E = "  -3.07%  " }
    35 = @{ D = "0.7518"; E = "  +0.19%  " }
    36 = @{ E = "  -2.11%  " }
    37 = @{ D = "2.710"; E = "  +0.46%  " }
    38 = @{ D = "0.01871"; E = "  -1.05%  " }
    39 = @{ D = "1.244.76"; E = "  -3.24%  " }
    40 = @{ D = "2.736"; E = "  -0.31%  " }
    41 = @{ D = "111.45"; E = "  -0.38%  " }
    42 = @{ D = "0.8990"; E = "  -1.69%  " }
    43 = @{ D = "6.125"; E = "  -7.14%  " }
    44 = @{ D = "70.91"; E = "  -4.98%  " }
    45 = @{ D = "0.9990"; E = "  -0.15%  " }
    46 = @{ E = "  -1.51%  " }
    47 = @{ D = "2.006.74"; E = "  -1.55%  " }
    48 = @{ D = "0.5190"; E = "  -0.57%  " }
    49 = @{ E = "  -1.43%  " }
    50 = @{ D = "9.503"; E = "  -0.37%  " }
    51 = @{ D = "0.4309"; E = "  -1.98%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
